$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1. Delete the "better meta models" task row entirely
for ($i = $tbl.ListRows.Count; $i -ge 1; $i--) {
    $row = $tbl.ListRows.Item($i)
    $val = $row.Range.Cells.Item(1,1).Value()
    if ($val -eq "better meta models") {
        $ws.Rows($row.Range.Row).Delete()
    }
}

# 2. Add the new "better lamps" task row
$newRow = $tbl.ListRows.Add()
$newRow.Range.Cells.Item(1,1).Value = "better lamps: add materials and swap material to a glowing one when on"
$newRow.Range.Cells.Item(1,2).Value = "QoL, Visual Only"
$newRow.Range.Cells.Item(1,3).Value = 80
$newRow.Range.Cells.Item(1,4).Value = 5

# 3. Re-sort the table by Priority then Workload (ascending), as before
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($ws.Range("C2:C" + (1 + $tbl.ListRows.Count)))
$tbl.Sort.SortFields.Add($ws.Range("D2:D" + (1 + $tbl.ListRows.Count)))
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# 4. Fix up the row heights for wrapped long-text rows so they match their (moved) content
$ws.Rows(6).RowHeight = 45
$ws.Rows(12).RowHeight = 30
$ws.Rows(13).RowHeight = 30
$ws.Rows(14).RowHeight = 30
$ws.Rows(16).RowHeight = 30
$ws.Rows(28).RowHeight = 30
$ws.Rows(31).RowHeight = 30
$ws.Rows(32).RowHeight = 30
$ws.Rows(33).RowHeight = 30
$ws.Rows(34).RowHeight = 30
$ws.Rows(35).RowHeight = 45

# 5. Update the selection to A2
$ws.Range("A2").Select()
